# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 74
$wsExhibit.Range("F3").Value = 11784
$wsExhibit.Range("F4").Value = 219
$wsExhibit.Range("F7").Value = 11722
$wsExhibit.Range("F10").Value = 96
$wsExhibit.Range("F11").Value = 36
$wsExhibit.Range("F12").Value = 1771
$wsExhibit.Range("F13").Value = 5807

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 74
$wsAll.Range("F5").Value = 11784
$wsAll.Range("F6").Value = 219
$wsAll.Range("F9").Value = 11722
$wsAll.Range("F12").Value = 96
$wsAll.Range("F13").Value = 36
$wsAll.Range("F14").Value = 1771
$wsAll.Range("F16").Value = 5807

$wb.Save()
